$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("condition_Library")

# Add new H-column (丧值变更 / san_ChangeVal) values of 0 for rows 3, 4, 6 and 7
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0

# Move the active selection in the bottom-right frozen pane to H7
$ws.Range("H7").Select()
